# Update the Markov transition-probability matrix on Sheet1 with the
# recomputed values after adding more simulated games (see commit message:
# "added more games, sped up simulate game logic, and drafted optimization
# logic"). Each affected row is a probability distribution that still sums
# to 1 after the update; only the cells whose underlying game counts
# changed are touched here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1930501930501931
$ws.Range("C2").Value = 0.5405405405405406
$ws.Range("J2").Value = 0.03088803088803089
$ws.Range("P2").Value = 0.1312741312741313
$ws.Range("S2").Value = 0.1042471042471042

$ws.Range("B3").Value = 0.02836879432624113
$ws.Range("C3").Value = 0.007092198581560284
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.7021276595744681
$ws.Range("S3").Value = 0.2411347517730496

$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.6551724137931034
$ws.Range("S4").Value = 0.2758620689655172

$ws.Range("B6").Value = 0.08888888888888889
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("F6").Value = 0.07222222222222222
$ws.Range("J6").Value = 0.1777777777777778
$ws.Range("O6").Value = 0.01111111111111111
$ws.Range("Q6").Value = 0.1277777777777778
$ws.Range("R6").Value = 0.1
$ws.Range("S6").Value = 0.4055555555555556

$ws.Range("B7").Value = 0.1023622047244094
$ws.Range("F7").Value = 0.03937007874015748
$ws.Range("J7").Value = 0.1338582677165354
$ws.Range("O7").Value = 0.02362204724409449
$ws.Range("Q7").Value = 0.1968503937007874
$ws.Range("R7").Value = 0.07874015748031496
$ws.Range("S7").Value = 0.4251968503937008

$ws.Range("B8").Value = 0.1143617021276596
$ws.Range("D8").Value = 0.01063829787234043
$ws.Range("E8").Value = 0.002659574468085106
$ws.Range("F8").Value = 0.05851063829787234
$ws.Range("J8").Value = 0.1329787234042553
$ws.Range("O8").Value = 0.02393617021276596
$ws.Range("Q8").Value = 0.1648936170212766
$ws.Range("R8").Value = 0.06648936170212766
$ws.Range("S8").Value = 0.425531914893617

$ws.Range("B9").Value = 0.09259259259259259
$ws.Range("D9").Value = 0.03703703703703703
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.1358024691358025
$ws.Range("O9").Value = 0.02469135802469136
$ws.Range("Q9").Value = 0.1111111111111111
$ws.Range("R9").Value = 0.06172839506172839
$ws.Range("S9").Value = 0.4814814814814815

$ws.Range("B10").Value = 0.1086556169429098
$ws.Range("D10").Value = 0.0147329650092081
$ws.Range("F10").Value = 0.07642725598526703
$ws.Range("J10").Value = 0.1261510128913444
$ws.Range("O10").Value = 0.01933701657458563
$ws.Range("Q10").Value = 0.1813996316758748
$ws.Range("R10").Value = 0.07734806629834254
$ws.Range("S10").Value = 0.3959484346224678

$ws.Range("G11").Value = 0.1778846153846154
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.2259615384615385
$ws.Range("L11").Value = 0.4903846153846154
$ws.Range("S11").Value = 0.02884615384615385

$ws.Range("G12").Value = 0.6981132075471698
$ws.Range("J12").Value = 0.2264150943396226
$ws.Range("L12").Value = 0.02830188679245283
$ws.Range("S12").Value = 0.04716981132075472

$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.03846153846153846

$ws.Range("H15").Value = 0.1741293532338309
$ws.Range("I15").Value = 0.0945273631840796
$ws.Range("J15").Value = 0.3880597014925373
$ws.Range("K15").Value = 0.02985074626865672
$ws.Range("M15").Value = 0.004975124378109453
$ws.Range("O15").Value = 0.05472636815920398
$ws.Range("S15").Value = 0.2537313432835821

$ws.Range("F16").Value = 0.02649006622516556
$ws.Range("H16").Value = 0.2119205298013245
$ws.Range("I16").Value = 0.0728476821192053
$ws.Range("J16").Value = 0.423841059602649
$ws.Range("K16").Value = 0.1059602649006623
$ws.Range("M16").Value = 0.006622516556291391
$ws.Range("O16").Value = 0.05298013245033113
$ws.Range("S16").Value = 0.09933774834437085

$ws.Range("F17").Value = 0.01238390092879257
$ws.Range("H17").Value = 0.1609907120743034
$ws.Range("I17").Value = 0.09597523219814241
$ws.Range("J17").Value = 0.4705882352941176
$ws.Range("K17").Value = 0.07739938080495357
$ws.Range("M17").Value = 0.02476780185758514
$ws.Range("O17").Value = 0.0804953560371517
$ws.Range("S17").Value = 0.07739938080495357

$ws.Range("F18").Value = 0.0273972602739726
$ws.Range("H18").Value = 0.2191780821917808
$ws.Range("I18").Value = 0.07534246575342465
$ws.Range("J18").Value = 0.4520547945205479
$ws.Range("K18").Value = 0.0410958904109589
$ws.Range("M18").Value = 0.0136986301369863
$ws.Range("O18").Value = 0.0684931506849315
$ws.Range("S18").Value = 0.1027397260273973

$ws.Range("F19").Value = 0.008571428571428572
$ws.Range("H19").Value = 0.2142857142857143
$ws.Range("I19").Value = 0.08761904761904762
$ws.Range("J19").Value = 0.4123809523809524
$ws.Range("K19").Value = 0.1009523809523809
$ws.Range("M19").Value = 0.01238095238095238
$ws.Range("O19").Value = 0.08190476190476191
$ws.Range("S19").Value = 0.08190476190476191
